# Slide 1's title placeholder text began with two separately-authored runs:
# "Pré" (flagged err="1" by the spell checker) and " projeto". The commit
# "Apresentação sem formatação pronta" drops the leading "Pré " so the title
# now simply starts with "projeto" - i.e. the two runs collapse into a
# single run whose text is "projeto".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$needle = "Pré projeto"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $idx = $tr.Text.IndexOf($needle)
        if ($idx -ge 0) {
            # Grab just the "Pré projeto" span and collapse it to "projeto",
            # leaving every other run (sizes, line breaks, later text) intact.
            $sub = $tr.Characters($idx + 1, $needle.Length)
            $sub.Text = "projeto"
            break
        }
    }
}
